# "Generate Report for Handoff" - refresh the localization-status report:
#  - Rows 4-7 (status "Ready for handoff") move from Priority "low" to "ht",
#    matching the already-handed-off rows 2-3 on both language sheets.
#  - The handoff timestamp for the 120c0da6 file (zh-cn) advances from
#    00:36:15 to 00:36:42.
#  - The shared "Latest HO Xliff Generate Date" timestamp used by the
#    Overview sheet and the de-de handoff rows advances from
#    00:36:23 to 00:36:50.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet (rows 4-7): Priority low -> ht, Latest Handoff Datetime refreshed
for ($r = 4; $r -le 7; $r++) {
    $wsZhCn.Cells.Item($r, 5).Value = "ht"
    $wsZhCn.Cells.Item($r, 8).Value = "2016-08-13 00:36:42"
}

# de-de sheet (rows 4-7): Priority low -> ht, Latest Handoff Datetime refreshed
for ($r = 4; $r -le 7; $r++) {
    $wsDeDe.Cells.Item($r, 5).Value = "ht"
    $wsDeDe.Cells.Item($r, 8).Value = "2016-08-13 00:36:50"
}

# Overview sheet (rows 4-7): Latest HO Xliff Generate Date refreshed
for ($r = 4; $r -le 7; $r++) {
    $wsOverview.Cells.Item($r, 7).Value = "2016-08-13 00:36:50"
}
